$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Applies the updated crypto price/volume snapshot cell by cell.
# For numeric-looking strings in column D we prefix with an apostrophe
# so Excel keeps them as literal text (preserving trailing zeros / the
# multi-dot "thousands.hundreds" notation used by the source feed)
# instead of silently coercing them to a Double.

$ws.Range("D2").Value = "24.911.02"
$ws.Range("E2").Value = "  -3.30%  "
$ws.Range("D3").Value = "1.675.04"
$ws.Range("E3").Value = "  -3.19%  "
$ws.Range("D4").Value = "'1.006"
$ws.Range("E4").Value = "  +0.62%  "
$ws.Range("D5").Value = "'309.55"
$ws.Range("E5").Value = "  -1.70%  "
$ws.Range("D6").Value = "'0.9994"
$ws.Range("D7").Value = "'0.3665"
$ws.Range("E7").Value = "  -3.19%  "
$ws.Range("D8").Value = "'0.3367"
$ws.Range("E8").Value = "  -7.10%  "
$ws.Range("D9").Value = "'47.70"
$ws.Range("E9").Value = "  -5.65%  "
$ws.Range("D10").Value = "'1.169"
$ws.Range("E10").Value = "  -4.12%  "
$ws.Range("D11").Value = "'0.07296"
$ws.Range("E11").Value = "  -4.44%  "
$ws.Range("D12").Value = "'0.9988"
$ws.Range("E12").Value = "  +0.36%  "
$ws.Range("D13").Value = "'6.162"
$ws.Range("E13").Value = "  -4.59%  "
$ws.Range("D14").Value = "'20.39"
$ws.Range("E14").Value = "  -5.88%  "
$ws.Range("D15").Value = "'6.796"
$ws.Range("E15").Value = "  -3.94%  "
$ws.Range("D16").Value = "1.674.06"
$ws.Range("E16").Value = "  -3.40%  "
$ws.Range("D17").Value = "'0.00001095"
$ws.Range("E17").Value = "  -4.80%  "
$ws.Range("D18").Value = "'0.9996"
$ws.Range("E18").Value = "  +0.70%  "
$ws.Range("D19").Value = "'0.06588"
$ws.Range("E19").Value = "  -3.18%  "
$ws.Range("D20").Value = "'81.86"
$ws.Range("E20").Value = "  -5.25%  "
$ws.Range("D21").Value = "'16.78"
$ws.Range("E21").Value = "  -3.42%  "
$ws.Range("D22").Value = "'6.143"
$ws.Range("E22").Value = "  -4.80%  "
$ws.Range("E23").Value = "  -0.87%  "
$ws.Range("D24").Value = "24.916.19"
$ws.Range("E24").Value = "  -3.09%  "
$ws.Range("D25").Value = "'2.436"
$ws.Range("E25").Value = "  -0.21%  "
$ws.Range("D26").Value = "'2.671"
$ws.Range("E26").Value = "  -9.27%  "
$ws.Range("D27").Value = "'19.76"
$ws.Range("E27").Value = "  -3.66%  "
$ws.Range("D28").Value = "'149.69"
$ws.Range("D29").Value = "'1.248"
$ws.Range("E29").Value = "  +4.15%  "
$ws.Range("D30").Value = "'130.13"
$ws.Range("E30").Value = "  -3.81%  "
$ws.Range("D31").Value = "1.857.60"
$ws.Range("E31").Value = "  -3.50%  "
$ws.Range("D32").Value = "'6.456"
$ws.Range("E32").Value = "  -6.14%  "
$ws.Range("D33").Value = "'4.152"
$ws.Range("E33").Value = "  +1.25%  "
$ws.Range("B34").Value = "WEMIXTOKEN"
$ws.Range("C34").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D34").Value = "'1.739"
$ws.Range("E34").Value = "  -3.49%  "
$ws.Range("B35").Value = "Aptos"
$ws.Range("C35").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D35").Value = "'13.36"
$ws.Range("E35").Value = "  -3.77%  "
$ws.Range("B36").Value = "Stellar"
$ws.Range("C36").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D36").Value = "'0.08583"
$ws.Range("E36").Value = "  -0.51%  "
$ws.Range("E37").Value = "  -3.75%  "
$ws.Range("D38").Value = "'0.06432"
$ws.Range("E38").Value = "  -4.59%  "
$ws.Range("E39").Value = "  -5.28%  "
$ws.Range("D40").Value = "'8.678"
$ws.Range("E40").Value = "  -6.15%  "
$ws.Range("D41").Value = "'0.2150"
$ws.Range("E41").Value = "  -2.84%  "
$ws.Range("D42").Value = "'1.248"
$ws.Range("E42").Value = "  -3.75%  "
$ws.Range("D43").Value = "'0.6246"
$ws.Range("E43").Value = "  -3.36%  "
$ws.Range("D44").Value = "'0.9985"
$ws.Range("E44").Value = "  +0.52%  "
$ws.Range("D45").Value = "'13.40"
$ws.Range("E45").Value = "  -2.81%  "
$ws.Range("D46").Value = "'3.793"
$ws.Range("E46").Value = "  -2.56%  "
$ws.Range("E47").Value = "  -5.05%  "
$ws.Range("D48").Value = "'2.033"
$ws.Range("E48").Value = "  -5.35%  "
$ws.Range("D49").Value = "'125.35"
$ws.Range("E49").Value = "  -4.61%  "
$ws.Range("D50").Value = "'0.07154"
$ws.Range("E50").Value = "  -3.84%  "
$ws.Range("D51").Value = "'76.97"
$ws.Range("E51").Value = "  -2.64%  "
